# Add pk_control column (I) with header and incrementing control values,
# and change column B's date cells to a custom yyyy-mm-dd number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "pk_control" header + values in column I
$ws.Range("I1").Value = "pk_control"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = $r + 42
}

# Apply custom date format (yyyy-mm-dd) to column B header + data
$ws.Range("B1").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B2:B21").NumberFormat = "yyyy\-mm\-dd"
$ws.Columns.Item(2).ColumnWidth = 10

# Reselect cell A21 on the active sheet (matches last worked-on cell)
[void]$ws.Range("A21").Select()
